$d = $word.ActiveDocument

# "Change Libfabric back to libfabric - spell check fix"
#
# The document has a mix of correctly-spelled "libfabric" and incorrectly
# capitalised "Libfabric" occurrences (plus a couple of places where
# "libfabric" is awkwardly split across two runs as "l" + "ibfabric").
# This script normalises all of them to the single, lower-case word
# "libfabric", using enough surrounding context in each Find so that only
# the intended occurrence is touched (several correctly-spelled
# "libfabric" instances elsewhere in the document, e.g. in the references
# section, must be left completely untouched).

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null
}

# 1) "Our new parcelport implementation is based on Libfabric " -> libfabric
Replace-Text "based on Libfabric " "based on libfabric "

# 2) "...to improve performance and Libfabric also supports..." -> libfabric
Replace-Text "performance and Libfabric also" "performance and libfabric also"

# 3) "...the speed improvement of the Libfabric parcelport..." -> libfabric
Replace-Text "improvement of the Libfabric parcelport" "improvement of the libfabric parcelport"

# 4) "...refinement (LoR) with the " + "l" + "ibfabric and MPI parcelports..."
#    -> merge the split "l"/"ibfabric" runs back into a single "libfabric" run
Replace-Text "(LoR) with the libfabric and MPI" "(LoR) with the libfabric and MPI"

# 5) "...between two nodes using the " + "l" + "ibfabric" + " parcelport in HPX..."
Replace-Text "using the libfabric parcelport in HPX on different thread counts" "using the libfabric parcelport in HPX on different thread counts"

# 6) "...would like to thank the " + "l" + "ibfabric" + " developers..."
Replace-Text "thank the libfabric developers" "thank the libfabric developers"
